# Re-sort the GAM diagnostics table (rows 2-9) into ascending order by
# model name (g1..g8). The underlying per-model data (covariates, aic,
# gcv, r_squared, dev_expl) travels with its row - only the row order
# changes, fixing the mis-ordered confidence-interval/summary table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:F9")
$sortKey = $ws.Range("A2:A9")

$dataRange.Sort($sortKey, 1)
